$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Shravya2"
$ws.Range("B4").Value = 3456

$ws.Range("B4").Select()
